$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, bordered, centered) used by the other
# header cells (e.g. H1) by copying H1's formatting onto the new headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for the new I (I0) and J (IF) columns, keyed by row number.
# I is 1 for almost every row (a couple of exceptions), and J mirrors the
# existing H (IP) value for almost every row (again, a couple of exceptions).
$colI = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 1; 27 = 1; 28 = 1;
    29 = 1; 30 = 1; 31 = 1; 32 = 1; 33 = 1; 34 = 1; 35 = 2; 36 = 1; 37 = 1;
    38 = 3; 39 = 1;
}
$colJ = @{
    2 = 5; 3 = 6; 4 = 5; 5 = 4; 6 = 4; 7 = 9; 8 = 5; 9 = 5; 10 = 5;
    11 = 7; 12 = 7; 13 = 6; 14 = 4; 15 = 7; 16 = 7; 17 = 6; 18 = 7; 19 = 6;
    20 = 5; 21 = 7; 22 = 7; 23 = 5; 24 = 6; 25 = 6; 26 = 6; 27 = 5; 28 = 3;
    29 = 6; 30 = 6; 31 = 7; 32 = 7; 33 = 5; 34 = 5; 35 = 5; 36 = 4; 37 = 3;
    38 = 5; 39 = 2;
}

for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 9).Value = $colI[$r]
    $ws.Cells.Item($r, 10).Value = $colJ[$r]
}
